$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows above the current row 645; this shifts the
# existing rows 645:721 down to 649:725 (and the dimension grows to
# A1:T725 automatically).
$ws.Rows.Item(645).EntireRow.Insert()
$ws.Rows.Item(645).EntireRow.Insert()
$ws.Rows.Item(645).EntireRow.Insert()
$ws.Rows.Item(645).EntireRow.Insert()

# Common column values shared by every data row in this block.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100106
$producto  = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"
$origen    = "Provincia de Limarí"
$kgUnidad  = 1
$unidad17  = "$/kilo (en caja de 17 kilos)"

# New rows of data (date serial, variedad, calidad, volumen, precio min,
# precio max, precio promedio ponderado).
$newRows = @(
    @{ Row = 645; Fecha = 45106; Variedad = "Hass";              Calidad = "Especial"; Volumen = 100; PMin = 4800; PMax = 4900; PProm = 4850 },
    @{ Row = 646; Fecha = 45106; Variedad = "Hass";              Calidad = "Primera";  Volumen = 200; PMin = 4500; PMax = 4600; PProm = 4550 },
    @{ Row = 647; Fecha = 45106; Variedad = "Hass";              Calidad = "Segunda";  Volumen = 120; PMin = 4200; PMax = 4300; PProm = 4250 },
    @{ Row = 648; Fecha = 45106; Variedad = "Negra de La Cruz";  Calidad = "Primera";  Volumen = 160; PMin = 1900; PMax = 2000; PProm = 1950 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad17
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PProm
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
